$p = $ppt.ActivePresentation

# Remove the slide "Text Mining: Term Frequency" (sldId 258), which is the
# 4th slide in the deck's show order.
$p.Slides.Item(4).Delete()
